# Insert a new data row for "Pomelo" (Start Ruby, Primera) at sheet row 201.
# This pushes the existing rows 201-302 down to 202-303 (Excel preserves
# their values/formatting automatically), and the sheet's used range grows
# from A1:T302 to A1:T303.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(201).Insert()

$ws.Cells.Item(201, 1).Value  = 10
$ws.Cells.Item(201, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(201, 3).Value  = "La Araucanía"
$ws.Cells.Item(201, 4).Value  = 44813
$ws.Cells.Item(201, 5).Value  = 9
$ws.Cells.Item(201, 6).Value  = "Fruta"
$ws.Cells.Item(201, 7).Value  = 100102
$ws.Cells.Item(201, 8).Value  = "Cítricos"
$ws.Cells.Item(201, 9).Value  = 100102006
$ws.Cells.Item(201, 10).Value = "Pomelo"
$ws.Cells.Item(201, 11).Value = "Start Ruby"
$ws.Cells.Item(201, 12).Value = "Primera"
$ws.Cells.Item(201, 13).Value = 75
$ws.Cells.Item(201, 14).Value = 13000
$ws.Cells.Item(201, 15).Value = 13000
$ws.Cells.Item(201, 16).Value = 13000
$ws.Cells.Item(201, 17).Value = "$/bandeja 15 kilos granel"
$ws.Cells.Item(201, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(201, 19).Value = 867
$ws.Cells.Item(201, 20).Value = 15
